$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-21 (Title, Genre, Rating, Year)
$data = @(
    @("Love Lies Bleeding", "Action", "7", "2024"),
    @("Arthur the King", "Adventure", "7", "2024"),
    @("[FR] Infested", "Horror", "6.7", "2023"),
    @("Monkey Man", "Action", "7.1", "2024"),
    @("The Old Oak", "Action", "7.1", "2023"),
    @("The Mitchells vs. the Machines", "Action", "7.6", "2021"),
    @("The Commuter", "Action", "6.3", "2018"),
    @("Justice League: Crisis on Infinite Earths - Part Two", "Action", "5.6", "2024"),
    @("Alice Through the Looking Glass", "Action", "6.2", "2016"),
    @("The Little Things", "Action", "6.3", "2021"),
    @("Good Burger 2", "Action", "5.3", "2023"),
    @("The Pod Generation", "Action", "5.6", "2023"),
    @("Gods of Egypt", "Action", "5.4", "2016"),
    @("Mean Girls", "Comedy", "5.7", "2024"),
    @("Drive-Away Dolls", "Action", "5.5", "2024"),
    @("Problemista", "Action", "7", "2023"),
    @("The Channel", "Action", "5.6", "2023"),
    @("Rebel Moon - Part Two: The Scargiver", "Action", "5.2", "2024"),
    @("Late Night with the Devil", "Horror", "7.1", "2023"),
    @("One Life", "Action", "7.6", "2023")
)

# Preserve text formatting for the Rating/Year columns so numeric-looking
# strings stay stored as text, matching the original inline-string cell
# typing rather than being coerced into numeric cells.
$ws.Range("C2:D21").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
